$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: kode "STF" (with name "fitria" originally referenced via B3) is
# replaced - the level_kode row now reads DIREKTUR / syafa.
$ws.Range("A2").Value = "DIREKTUR"
$ws.Range("B2").Value = "syafa"

# Remove the old rows 3 and 4 (STF/fitria and STF/hertin) entirely so the
# table only has the header plus the single DIREKTUR/syafa row left.
$ws.Rows("3:4").EntireRow.Delete() | Out-Null

$ws.Range("B2").Select() | Out-Null
